$d = $word.ActiveDocument

# --- Step 1: split paragraph 3 ("Na de eerste weken ... CONST.") ---
# Insert a new paragraph right after paragraph 3; the bookmark _GoBack that
# used to live inside paragraph 3 (before "loop weer onder de knie") will be
# re-created later inside the NEW paragraph, at the very end of its text.
$p3 = $d.Paragraphs(3)
$r = $p3.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# --- Step 2: fill the new paragraph (#4) with its text ---
$p4 = $d.Paragraphs(4)
# Temporarily append marker characters so the bookmark insertion point is not
# located exactly at the paragraph end (a boundary case that this runtime
# mishandles), then trim the marker off again.
$p4.Range.Text = "Het moeilijkste wat ik moest doen was werken met Objects maar dat werkt na tijd duidlijke.XX"
$p4b = $d.Paragraphs(4)
$bmPos = $p4b.Range.End - 3
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$trimStart = $p4b.Range.End - 3
$trimRange = $d.Range($trimStart, $trimStart + 2)
$trimRange.Text = ""

# --- Step 3: insert the remaining new paragraphs after paragraph 4 ---
$p4c = $d.Paragraphs(4)
$ins = $p4c.Range
$ins.Collapse(0)
$ins.InsertParagraphAfter()

# paragraph 5: empty
$p5 = $d.Paragraphs(5)
$ins2 = $p5.Range
$ins2.Collapse(0)
$ins2.InsertParagraphAfter()

# paragraph 6: "Raaf van Hulst ..."
$p6 = $d.Paragraphs(6)
$p6.Range.Text = "Raaf van Hulst en Tim van Eert hebben samen mijn stemwijzer na gekeken en vonden het volgende."

$p6b = $d.Paragraphs(6)
$ins3 = $p6b.Range
$ins3.Collapse(0)
$ins3.InsertParagraphAfter()

# paragraph 7: empty
$p7 = $d.Paragraphs(7)
$ins4 = $p7.Range
$ins4.Collapse(0)
$ins4.InsertParagraphAfter()

# paragraph 8: "Dat ik duidelijker benamingen ..."
$p8 = $d.Paragraphs(8)
$p8.Range.Text = "Dat ik duidelijker benamingen moet gaan gebruiken daarnaast willen ze dat ik ook meer CONSTS gaat gebruiken."

$p8b = $d.Paragraphs(8)
$ins5 = $p8b.Range
$ins5.Collapse(0)
$ins5.InsertParagraphAfter()

# paragraph 9: "Buiten dat vinden ze dat alle beoordelingscriteria voldoende zijn."
$p9 = $d.Paragraphs(9)
$p9.Range.Text = "Buiten dat vinden ze dat alle beoordelingscriteria voldoende zijn."

$p9b = $d.Paragraphs(9)
$ins6 = $p9b.Range
$ins6.Collapse(0)
$ins6.InsertParagraphAfter()

# paragraph 10: "De filter van partijen en of ze wel of niet seculair zijn."
$p10 = $d.Paragraphs(10)
$p10.Range.Text = "De filter van partijen en of ze wel of niet seculair zijn."
